# Commit: "changed data function to include cr-28 data"
#
# Insert a new worksheet "cr-28 var" between "Problem Variables" and
# "Dynamometer", populated with a CR-28 variant of the "Problem Variables"
# data (same layout, a handful of updated parameter values, and a single
# "CR-28" row in the MODEL table instead of the five comparison cars).

$wb = $excel.ActiveWorkbook

$problemVars = $wb.Worksheets.Item("Problem Variables")
$dyno = $wb.Worksheets.Item("Dynamometer")

# New sheet goes right after "Problem Variables" (i.e. right before "Dynamometer").
$ws = $wb.Worksheets.Add($null, $problemVars)
$ws.Name = "cr-28 var"

# ---- Title -----------------------------------------------------------
$ws.Range("A1").Value = "DATA FOR THE CAR PERFORMANCE CR-28"

# ---- GEAR RATIOS -------------------------------------------------------
$ws.Range("A3").Value = "GEAR RATIOS"
$ws.Range("B3").Value = "value"
$ws.Range("C3").Value = "test"

$ws.Range("A4").Value = "1st gear"
$ws.Range("B4").Value = 4.171
$ws.Range("C4").Value = 4.171

$ws.Range("A5").Value = "2nd gear"
$ws.Range("B5").Value = 2.34
$ws.Range("C5").Value = 2.34

$ws.Range("A6").Value = "3rd gear"
$ws.Range("B6").Value = 1.521
$ws.Range("C6").Value = 1.521

$ws.Range("A7").Value = "4th gear"
$ws.Range("B7").Value = 1.143
$ws.Range("C7").Value = 1.143

$ws.Range("A8").Value = "5th gear"
$ws.Range("B8").Value = 0.867
$ws.Range("C8").Value = 0.867

$ws.Range("A9").Value = "6th gear"
$ws.Range("B9").Value = 0.691
$ws.Range("C9").Value = 0.691

# ---- OTHER PARAMETERS ---------------------------------------------------
$ws.Range("A11").Value = "OTHER PARAMETERS"
$ws.Range("B11").Value = "symbol"
$ws.Range("C11").Value = "value"
$ws.Range("D11").Value = "test"

$ws.Range("A12").Value = "speed in km/h"
$ws.Range("B12").Value = "v"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 100

$ws.Range("A13").Value = "terrain slope as %"
$ws.Range("B13").Value = "slope %"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 2

$ws.Range("A14").Value = "wheelbase in m"
$ws.Range("B14").Value = "L"
$ws.Range("C14").Value = 1.53
$ws.Range("D14").Value = 1.53

$ws.Range("A15").Value = "radius of the wheel in m"
$ws.Range("B15").Value = "r"
$ws.Range("C15").Value = 0.2
$ws.Range("D15").Value = 0.2

$ws.Range("A16").Value = "roll resistance coefficient"
$ws.Range("B16").Value = "f"
$ws.Range("C16").Value = 0.012
$ws.Range("D16").Value = 0.012

$ws.Range("A17").Value = "air resistance moment arm / center of gravity height in m"
$ws.Range("B17").Value = "hA/h"
$ws.Range("C17").Value = 0.5
$ws.Range("D17").Value = 0.5

$ws.Range("A18").Value = "final drive efficiency as %"
$ws.Range("B18").Value = "ηd"
$ws.Range("C18").Value = 95
$ws.Range("D18").Value = 95

$ws.Range("A19").Value = "transmission efficiency as %"
$ws.Range("B19").Value = "ηt"
$ws.Range("C19").Value = 98
$ws.Range("D19").Value = 98

$ws.Range("A20").Value = "weight of the vehicle in N"
$ws.Range("B20").Value = "W"
$ws.Range("C20").Value = 2134
$ws.Range("D20").Value = 2134

$ws.Range("A21").Value = "air density in kg/m^3"
$ws.Range("B21").Value = "ρ"
$ws.Range("C21").Value = 1.225
$ws.Range("D21").Value = 1.225

$ws.Range("A22").Value = "drive ratio"
$ws.Range("B22").Value = "id"
$ws.Range("C22").Value = 3.5
$ws.Range("D22").Value = 3.5

$ws.Range("A23").Value = "distance identifying the center of gravity position from the front wheel in m"
$ws.Range("B23").Value = "b"
$ws.Range("C23").Value = 0.765
$ws.Range("D23").Value = 0.765

# ---- MODEL ---------------------------------------------------------------
$ws.Range("A25").Value = "MODEL"
$ws.Range("B25").Value = "Year"
$ws.Range("C25").Value = "CD"
$ws.Range("D25").Value = "A (m^2)"

$ws.Range("A26").Value = "CR-28"
$ws.Range("B26").Value = 2023
$ws.Range("C26").Value = 1.81
$ws.Range("D26").Value = 1.07

# ---- Selection / active sheet --------------------------------------------
$problemVars.Activate() | Out-Null
$problemVars.Range("A1:E30").Select() | Out-Null

$ws.Activate() | Out-Null
$ws.Range("D24").Select() | Out-Null
